# Incontactables.xlsx — "Add files via upload" re-save.
#
# The committed diff corresponds to deleting one data row from Hoja1:
#   row 441 -> Delv Ship-To Party "4000308422" / "Chedraui 169 Xoxocotlan"
#              (Contactabilidad1 "Incontactable", Origen de cuenta "Chedraui")
# which sits right after a similarly-named "Dimuflo SA de CV" row and looks
# like a duplicate/erroneous entry that got cleaned up. Deleting the entire
# row shifts every subsequent row up by one (so the table now ends at row
# 894 instead of 895) and Excel garbage-collects the two shared strings
# ("4000308422" and "Chedraui 169 Xoxocotlan") that were only referenced by
# that row.
#
# The save also picked up an AutoFilter over the whole table (with its
# companion hidden _xlnm._FilterDatabase defined name) and the window/
# selection ended up scrolled down to around row 871 with B876 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the stray "4000308422 / Chedraui 169 Xoxocotlan" row ---------
$targetRow = 441
$idCheck = $ws.Range("A" + $targetRow).Value2
if ($idCheck -eq "4000308422") {
    $ws.Rows($targetRow).Delete()
} else {
    # Fallback: scan column A for the account id if the sheet shape ever
    # differs from what's expected, so the right row still gets removed.
    $lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 1).Value2 -eq "4000308422") {
            $ws.Rows($r).Delete()
            break
        }
    }
}

# --- Recompute the used range and (re)apply the AutoFilter ---------------
$lastRow = $ws.Cells($ws.Rows.Count, 1).End(-4162).Row
$lastCol = $ws.Cells(1, $ws.Columns.Count).End(-4159).Column
$dataRange = $ws.Range($ws.Cells(1, 1), $ws.Cells($lastRow, $lastCol))
$dataRange.AutoFilter()

# AutoFilter in real Excel always creates/updates a hidden, sheet-scoped
# "_xlnm._FilterDatabase" defined name pointing at the filtered range —
# recreate that here.
foreach ($nm in $ws.Names) {
    if ($nm.Name -like "*_FilterDatabase") {
        $nm.Delete()
    }
}
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=" + $ws.Name + "!" + $dataRange.Address(1, 1, 1, $false))
$filterName.Visible = $false

# --- Restore the saved selection / scroll position ------------------------
$ws.Range("B876").Select()
